$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 7411462.5
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 7411462.5
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 22234387.5
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -22234723.5
$ws.Range("H18").Value = 2874.5
$ws.Range("I18").Value = 2874.5
$ws.Range("K18").Value = 2874.5
$ws.Range("M18").Value = -2590.5
$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()
$ws.Range("H43").Value = 17548938
$ws.Range("J43").Value = 6133.4287
$ws.Range("L43").Value = 6133.4287
$ws.Range("N43").Value = -6271.4287
$ws.Range("H62").Value = 499.5
$ws.Range("I62").Value = 499.5
$ws.Range("K62").Value = 499.5
$ws.Range("M62").Value = 124.5
$ws.Range("H65").Value = 499.5
$ws.Range("I65").Value = 499.5
$ws.Range("K65").Value = 2497.5
$ws.Range("M65").Value = 622.5
$ws.Range("H70").Value = 335332.66
$ws.Range("J70").Value = 2999.5
$ws.Range("L70").Value = 8998.5
$ws.Range("N70").Value = -9538.5
$ws.Range("H73").Value = 335332.66
$ws.Range("J73").Value = 2999.5
$ws.Range("L73").Value = 8998.5
$ws.Range("N73").Value = -10870.5
$ws.Range("H107").Value = 2188.3333
$ws.Range("I107").Value = 2368.7058
$ws.Range("J107").Value = 1750.2858
$ws.Range("K107").Value = 2368.7058
$ws.Range("L107").Value = 1750.2858
$ws.Range("M107").Value = -448.7058000000002
$ws.Range("N107").Value = -5590.2858
$ws.Range("H132").Value = 2014.2587
$ws.Range("I132").Value = 1960.5272
$ws.Range("K132").Value = 5881.5816
$ws.Range("M132").Value = -3351.5816
$ws.Range("H135").Value = 1918.129
$ws.Range("I135").Value = 1663.9231
$ws.Range("J135").Value = 3240
$ws.Range("K135").Value = 14975.3079
$ws.Range("L135").Value = 29160
$ws.Range("M135").Value = -12440.3079
$ws.Range("N135").Value = -34230
$ws.Range("H137").Value = 1853371.6
$ws.Range("I137").Value = 1025.4
$ws.Range("K137").Value = 3076.2
$ws.Range("M137").Value = -526.2000000000003
$ws.Range("H138").Value = 7873.884
$ws.Range("J138").Value = 6307.952
$ws.Range("L138").Value = 18923.856
$ws.Range("N138").Value = -29203.856

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 163069.28
$ws.Range("I32").Value = 191789.1
$ws.Range("J32").Value = 10854.3
$ws.Range("K32").Value = 191789.1
$ws.Range("L32").Value = 10854.3
$ws.Range("M32").Value = -191502.1
$ws.Range("N32").Value = -11428.3
$ws.Range("H102").Value = 1731.2
$ws.Range("I102").Value = 1292.4584
$ws.Range("K102").Value = 1292.4584
$ws.Range("M102").Value = 329.5416
$ws.Range("H122").Value = 2527.4614
$ws.Range("I122").Value = 2133.9524
$ws.Range("K122").Value = 6401.8572
$ws.Range("M122").Value = -3951.8572
$ws.Range("H132").Value = 3992.8857
$ws.Range("I132").Value = 3132.2917
$ws.Range("J132").Value = 5870.5454
$ws.Range("K132").Value = 9396.875100000001
$ws.Range("L132").Value = 17611.6362
$ws.Range("M132").Value = -6866.875100000001
$ws.Range("N132").Value = -22671.6362
$ws.Range("H135").Value = 40249
$ws.Range("J135").Value = 40249
$ws.Range("L135").Value = 40249
$ws.Range("N135").Value = -50389

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3279.9524
$ws.Range("I86").Value = 1848.4166
$ws.Range("J86").Value = 5188.6665
$ws.Range("K86").Value = 1848.4166
$ws.Range("L86").Value = 5188.6665
$ws.Range("M86").Value = -725.4166
$ws.Range("N86").Value = -7434.6665
$ws.Range("H89").Value = 3279.9524
$ws.Range("I89").Value = 1848.4166
$ws.Range("J89").Value = 5188.6665
$ws.Range("K89").Value = 9242.083000000001
$ws.Range("L89").Value = 25943.3325
$ws.Range("M89").Value = -3626.083000000001
$ws.Range("N89").Value = -37175.3325
$ws.Range("H134").Value = 19568040
$ws.Range("I134").Value = 2289.0688
$ws.Range("J134").Value = 52944908
$ws.Range("K134").Value = 6867.2064
$ws.Range("L134").Value = 158834724
$ws.Range("M134").Value = -4332.2064
$ws.Range("N134").Value = -158839794

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4456.0454
$ws.Range("I31").Value = 2510.8333
$ws.Range("J31").Value = 4763.184
$ws.Range("K31").Value = 2510.8333
$ws.Range("L31").Value = 4763.184
$ws.Range("M31").Value = -2215.8333
$ws.Range("N31").Value = -5353.184
$ws.Range("H34").Value = 4456.0454
$ws.Range("I34").Value = 2510.8333
$ws.Range("J34").Value = 4763.184
$ws.Range("K34").Value = 2510.8333
$ws.Range("L34").Value = 4763.184
$ws.Range("M34").Value = -2308.8333
$ws.Range("N34").Value = -5167.184
$ws.Range("H99").Value = 2252532
$ws.Range("I99").Value = 6667850
$ws.Range("J99").Value = 44873.168
$ws.Range("K99").Value = 6667850
$ws.Range("L99").Value = 44873.168
$ws.Range("M99").Value = -6666352
$ws.Range("N99").Value = -47869.168
$ws.Range("H126").Value = 2252532
$ws.Range("I126").Value = 6667850
$ws.Range("J126").Value = 44873.168
$ws.Range("K126").Value = 20003550
$ws.Range("L126").Value = 134619.504
$ws.Range("M126").Value = -20001080
$ws.Range("N126").Value = -139559.504
$ws.Range("H134").Value = 3328.5789
$ws.Range("I134").Value = 2742.7
$ws.Range("J134").Value = 3979.5557
$ws.Range("K134").Value = 8228.099999999999
$ws.Range("L134").Value = 11938.6671
$ws.Range("M134").Value = -5693.099999999999
$ws.Range("N134").Value = -17008.6671

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 125033.875
$ws.Range("J23").Value = 166698.5
$ws.Range("L23").Value = 500095.5
$ws.Range("N23").Value = -500565.5
$ws.Range("H129").Value = 24248352
$ws.Range("I129").Value = 3946.375
$ws.Range("J129").Value = 38102300
$ws.Range("K129").Value = 11839.125
$ws.Range("L129").Value = 114306900
$ws.Range("M129").Value = -6839.125
$ws.Range("N129").Value = -114316900

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 15152710
$ws.Range("J102").Value = 2187.5
$ws.Range("L102").Value = 2187.5
$ws.Range("N102").Value = -5431.5
$ws.Range("H107").Value = 239178.89
$ws.Range("J107").Value = 212.75
$ws.Range("L107").Value = 212.75
$ws.Range("N107").Value = -4052.75
$ws.Range("H113").Value = 2905.0715
$ws.Range("J113").Value = 3734.0908
$ws.Range("L113").Value = 3734.0908
$ws.Range("N113").Value = -8074.0908
$ws.Range("H126").Value = 3001.3333
$ws.Range("I126").Value = 3001.3333
$ws.Range("K126").Value = 9003.999899999999
$ws.Range("M126").Value = -6533.999899999999
$ws.Range("H132").Value = 6726342.5
$ws.Range("I132").Value = 3324.3704
$ws.Range("J132").Value = 19692162
$ws.Range("K132").Value = 9973.111199999999
$ws.Range("L132").Value = 59076486
$ws.Range("M132").Value = -7443.111199999999
$ws.Range("N132").Value = -59081546

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1014.2
$ws.Range("I16").Value = 886.6667
$ws.Range("K16").Value = 886.6667
$ws.Range("M16").Value = -716.6667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9205.5
$ws.Range("I62").Value = 8734
$ws.Range("J62").Value = 9441.25
$ws.Range("K62").Value = 8734
$ws.Range("L62").Value = 9441.25
$ws.Range("M62").Value = -8110
$ws.Range("N62").Value = -10689.25
$ws.Range("H65").Value = 9205.5
$ws.Range("I65").Value = 8734
$ws.Range("J65").Value = 9441.25
$ws.Range("K65").Value = 43670
$ws.Range("L65").Value = 47206.25
$ws.Range("M65").Value = -40550
$ws.Range("N65").Value = -53446.25
$ws.Range("H132").Value = 100880.2
$ws.Range("I132").Value = 100880.2
$ws.Range("K132").Value = 302640.6
$ws.Range("M132").Value = -300110.6
$ws.Range("H137").Value = 129990
$ws.Range("J137").Value = 129990
$ws.Range("L137").Value = 129990
$ws.Range("N137").Value = -140190
